# Refresh the "Price" (D) and "Volume(1h)" (E) columns of the cryptos sheet
# with the latest scraped values (GitHub Actions nightly update).
#
# Note: several Price values (e.g. "1.004", "310.67") look numeric, and a
# plain .Value assignment would have Excel auto-convert them to real numbers
# (losing formatting such as trailing zeros). To keep them as literal text -
# matching how the sheet stores every Price/Volume cell - we prefix those
# with a leading apostrophe (Excel's "store as text" quote-prefix marker)
# and then reset the cell style to "Normal" so no stray text-format style
# lingers on the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.904.86'
$ws.Range('E2').Value = '  +0.61%  '
$ws.Range('D3').Value = '1.810.67'
$ws.Range('E3').Value = '  +1.72%  '
$ws.Range('D4').Value = '''1.004'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('D5').Value = '''310.67'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.07%  '
$ws.Range('D6').Value = '''1.003'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.05%  '
$ws.Range('D7').Value = '''0.4979'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -2.72%  '
$ws.Range('D8').Value = '''0.3918'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +3.39%  '
$ws.Range('D9').Value = '''0.09546'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +22.69%  '
$ws.Range('E10').Value = '  +1.40%  '
$ws.Range('D11').Value = '''40.95'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.57%  '
$ws.Range('D12').Value = '''6.412'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +3.44%  '
$ws.Range('D13').Value = '''1.003'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +0.04%  '
$ws.Range('D14').Value = '''20.42'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +1.63%  '
$ws.Range('D15').Value = '1.808.57'
$ws.Range('E15').Value = '  +1.92%  '
$ws.Range('D16').Value = '''7.269'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +1.43%  '
$ws.Range('D17').Value = '''0.00001122'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +4.63%  '
$ws.Range('D18').Value = '''92.10'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +0.90%  '
$ws.Range('D19').Value = '''0.06658'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +1.56%  '
$ws.Range('E20').Value = '  +0.05%  '
$ws.Range('D21').Value = '''17.12'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.69%  '
$ws.Range('D22').Value = '''5.916'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.06%  '
$ws.Range('D23').Value = '27.968.23'
$ws.Range('E23').Value = '  +0.65%  '
$ws.Range('E24').Value = '  +1.21%  '
$ws.Range('D25').Value = '''2.251'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.61%  '
$ws.Range('D26').Value = '''158.97'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.25%  '
$ws.Range('D27').Value = '2.018.11'
$ws.Range('E28').Value = '  +1.64%  '
$ws.Range('D29').Value = '''2.376'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +1.11%  '
$ws.Range('D30').Value = '''127.39'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +2.37%  '
$ws.Range('E31').Value = '  -0.48%  '
$ws.Range('D32').Value = '''1.031'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +0.17%  '
$ws.Range('D33').Value = '''5.557'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +1.52%  '
$ws.Range('D34').Value = '''3.615'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -0.46%  '
$ws.Range('D35').Value = '''0.06715'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -4.83%  '
$ws.Range('D36').Value = '''8.910'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +2.05%  '
$ws.Range('E37').Value = '  +0.56%  '
$ws.Range('E38').Value = '  +0.61%  '
$ws.Range('D39').Value = '''4.928'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -1.85%  '
$ws.Range('D40').Value = '''11.19'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -2.54%  '
$ws.Range('D41').Value = '''0.6155'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +1.29%  '
$ws.Range('E42').Value = '  +0.09%  '
$ws.Range('E43').Value = '  -0.59%  '
$ws.Range('D44').Value = '''13.11'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +0.52%  '
$ws.Range('D45').Value = '''1.292'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -2.14%  '
$ws.Range('E46').Value = '  -0.89%  '
$ws.Range('D47').Value = '''3.699'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.36%  '
$ws.Range('D48').Value = '''123.01'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -3.57%  '
$ws.Range('D49').Value = '''1.927'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +1.84%  '
$ws.Range('D50').Value = '''1.178'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -1.37%  '
$ws.Range('D51').Value = '''0.06758'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.79%  '
